function Set-TextValue($range, $value) {
    # Force the cell to hold the given value as literal text, matching the
    # original inlineStr cells (avoids Excel auto-converting numeric-looking
    # strings like "582.07" or "62.934.26" into floating point numbers, and
    # restores the default cell style afterwards so no spurious style diff
    # is introduced).
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "62.934.26"
Set-TextValue $ws.Range("E2") "  +0.02%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.594.45"
Set-TextValue $ws.Range("E3") "  +2.03%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Range("E4") "  -0.05%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "582.07"
Set-TextValue $ws.Range("E5") "  +2.30%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "147.62"
Set-TextValue $ws.Range("E6") "  +1.38%  "

# Row 7 - USDC
Set-TextValue $ws.Range("E7") "  -0.03%  "

# Row 8 - XRP
Set-TextValue $ws.Range("D8") "0.597"
Set-TextValue $ws.Range("E8") "  +2.47%  "

# Row 9 - Dogecoin
Set-TextValue $ws.Range("D9") "0.108"
Set-TextValue $ws.Range("E9") "  +3.17%  "

# Row 10 - Toncoin
Set-TextValue $ws.Range("D10") "5.66"
Set-TextValue $ws.Range("E10") "  +2.69%  "

# Row 11 - TRON
Set-TextValue $ws.Range("E11") "  +0.05%  "

# Row 12 - Cardano
Set-TextValue $ws.Range("D12") "0.354"
Set-TextValue $ws.Range("E12") "  +0.70%  "

# Row 13 - Avalanche
Set-TextValue $ws.Range("D13") "27.15"
Set-TextValue $ws.Range("E13") "  -0.30%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D14") "3.059.59"
Set-TextValue $ws.Range("E14") "  +2.07%  "

# Row 15 - WrappedBTC
Set-TextValue $ws.Range("D15") "62.784.94"
Set-TextValue $ws.Range("E15") "  -0.12%  "

# Row 16 - ShibaInu
Set-TextValue $ws.Range("D16") "0.0000147"
Set-TextValue $ws.Range("E16") "  +3.14%  "

# Row 17 - WrappedEther
Set-TextValue $ws.Range("D17") "2.603.08"
Set-TextValue $ws.Range("E17") "  +2.15%  "

# Row 18 - Chainlink
Set-TextValue $ws.Range("D18") "11.34"
Set-TextValue $ws.Range("E18") "  +0.34%  "

# Row 19 - BitcoinCash
Set-TextValue $ws.Range("D19") "342.60"
Set-TextValue $ws.Range("E19") "  +2.60%  "

# Row 20 - Polkadot
Set-TextValue $ws.Range("D20") "4.40"
Set-TextValue $ws.Range("E20") "  +1.75%  "

# Row 21 - Uniswap
Set-TextValue $ws.Range("D21") "6.77"
Set-TextValue $ws.Range("E21") "  -0.08%  "

# Row 22 - Dai
Set-TextValue $ws.Range("D22") "1.00"
Set-TextValue $ws.Range("E22") "  +0.08%  "

# Row 23 - LEO
Set-TextValue $ws.Range("E23") "  -1.50%  "

# Row 24 - Litecoin
Set-TextValue $ws.Range("D24") "66.76"
Set-TextValue $ws.Range("E24") "  +2.60%  "

# Row 25 - WrappedeETH
Set-TextValue $ws.Range("D25") "2.724.84"
Set-TextValue $ws.Range("E25") "  +2.42%  "

# Row 26 - Kaspa
Set-TextValue $ws.Range("D26") "0.169"
Set-TextValue $ws.Range("E26") "  +0.12%  "

# Row 27 - Fetch.AI
Set-TextValue $ws.Range("E27") "  +0.34%  "

# Rows 28-30 reordered: Binance-PegBSC-USD, Aptos, InternetComputer(DFINITY)
# Row 28 -> Binance-PegBSC-USD
Set-TextValue $ws.Range("B28") "Binance-PegBSC-USD"
Set-TextValue $ws.Range("C28") "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue $ws.Range("D28") "1.01"
Set-TextValue $ws.Range("E28") "  +0.61%  "

# Row 29 -> Aptos
Set-TextValue $ws.Range("B29") "Aptos"
Set-TextValue $ws.Range("C29") "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D29") "7.94"
Set-TextValue $ws.Range("E29") "  +9.01%  "

# Row 30 -> InternetComputer(DFINITY)
Set-TextValue $ws.Range("B30") "InternetComputer(DFINITY)"
Set-TextValue $ws.Range("C30") "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D30") "8.38"
Set-TextValue $ws.Range("E30") "  +0.43%  "

# Row 31 - SuiNetwork
Set-TextValue $ws.Range("D31") "1.45"
Set-TextValue $ws.Range("E31") "  -1.94%  "

# Row 32 - PancakeSwap
Set-TextValue $ws.Range("E32") "  +4.23%  "

# Row 33 - PEPE
Set-TextValue $ws.Range("D33") "0.0₃0821"
Set-TextValue $ws.Range("E33") "  +1.11%  "

# Row 34 - Bittensor
Set-TextValue $ws.Range("D34") "460.19"
Set-TextValue $ws.Range("E34") "  +12.71%  "

# Row 35 - Monero
Set-TextValue $ws.Range("D35") "176.60"
Set-TextValue $ws.Range("E35") "  +0.77%  "

# Row 36 - ImmutableX
Set-TextValue $ws.Range("D36") "1.61"
Set-TextValue $ws.Range("E36") "  +3.79%  "

# Rows 37-38 reordered: FirstDigitalUSD, PolygonEcosystemToken
# Row 37 -> FirstDigitalUSD
Set-TextValue $ws.Range("B37") "FirstDigitalUSD"
Set-TextValue $ws.Range("C37") "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D37") "1.00"
Set-TextValue $ws.Range("E37") "  +0.00%  "

# Row 38 -> PolygonEcosystemToken
Set-TextValue $ws.Range("B38") "PolygonEcosystemToken"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue $ws.Range("D38") "0.404"
Set-TextValue $ws.Range("E38") "  +1.41%  "

# Row 39 - EthereumClassic
Set-TextValue $ws.Range("D39") "19.19"
Set-TextValue $ws.Range("E39") "  +1.20%  "

# Row 40 - NEARProtocol
Set-TextValue $ws.Range("D40") "4.54"
Set-TextValue $ws.Range("E40") "  +4.71%  "

# Row 41 - USDe
Set-TextValue $ws.Range("E41") "  +0.01%  "

# Row 42 - Stacks
Set-TextValue $ws.Range("D42") "1.70"
Set-TextValue $ws.Range("E42") "  -2.72%  "

# Row 43 - Aave
Set-TextValue $ws.Range("D43") "160.49"
Set-TextValue $ws.Range("E43") "  +5.70%  "

# Row 44 - Filecoin
Set-TextValue $ws.Range("D44") "3.79"
Set-TextValue $ws.Range("E44") "  +0.99%  "

# Row 45 - Mantle
Set-TextValue $ws.Range("D45") "0.636"
Set-TextValue $ws.Range("E45") "  +5.77%  "

# Row 46 - InjectiveProtocol
Set-TextValue $ws.Range("D46") "20.66"
Set-TextValue $ws.Range("E46") "  -0.36%  "

# Row 47 - Hedera
Set-TextValue $ws.Range("D47") "0.0545"
Set-TextValue $ws.Range("E47") "  +2.96%  "

# Row 48 - Stellar
Set-TextValue $ws.Range("D48") "0.0972"
Set-TextValue $ws.Range("E48") "  +1.06%  "

# Row 49 - VeChain
Set-TextValue $ws.Range("D49") "0.0237"
Set-TextValue $ws.Range("E49") "  -0.38%  "

# Row 50 - EnergySwap
Set-TextValue $ws.Range("D50") "18.55"
Set-TextValue $ws.Range("E50") "  +1.61%  "

# Row 51 - dogwifhat
Set-TextValue $ws.Range("D51") "1.73"
Set-TextValue $ws.Range("E51") "  -0.07%  "
